# Update data: 2025-10-29 10:19
# Refreshes the "Top Gainers", "Top Losers" and "1 Month Performance" tables
# with the latest market snapshot. Rows are overwritten in place (ticker +
# metric columns) to reflect the new rankings/values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Top Gainers")
$ws.Range("B55").Value = "JKIL"
$ws.Range("C55").Value = 4.1372
$ws.Range("D55").Value = 2.9463
$ws.Range("E55").Value = 1.7584
$ws.Range("B56").Value = "SAMBHV"
$ws.Range("C56").Value = 4.1349
$ws.Range("D56").Value = 2.624
$ws.Range("E56").Value = 5.167
$ws.Range("B57").Value = "PVRINOX"
$ws.Range("C57").Value = 4.1118
$ws.Range("D57").Value = 6.2102
$ws.Range("E57").Value = 14.707
$ws.Range("B58").Value = "KERNEX"
$ws.Range("C58").Value = 4.0782
$ws.Range("D58").Value = 7.542
$ws.Range("E58").Value = 27.2033
$ws.Range("B59").Value = "SUNFLAG"
$ws.Range("C59").Value = 3.997
$ws.Range("D59").Value = 4.333
$ws.Range("E59").Value = 4.6312
$ws.Range("C60").Value = 3.9096
$ws.Range("D60").Value = 2.6872
$ws.Range("E60").Value = 2.8935
$ws.Range("B64").Value = "ORIENTTECH"
$ws.Range("C64").Value = 3.827
$ws.Range("D64").Value = 0.5247
$ws.Range("E64").Value = 32.6784
$ws.Range("B65").Value = "SALASAR"
$ws.Range("C65").Value = 3.7935
$ws.Range("D65").Value = 4.7872
$ws.Range("E65").Value = 11.0485
$ws.Range("B66").Value = "NPST"
$ws.Range("C66").Value = 3.7841
$ws.Range("D66").Value = -2.0689
$ws.Range("E66").Value = -3.5677
$ws.Range("B67").Value = "DCW"
$ws.Range("C67").Value = 3.7544
$ws.Range("D67").Value = 2.3219
$ws.Range("E67").Value = -3.9753
$ws.Range("B71").Value = "BHARTIHEXA"
$ws.Range("C71").Value = 3.6718
$ws.Range("D71").Value = 7.0877
$ws.Range("E71").Value = 15.3332
$ws.Range("B72").Value = "HLEGLAS"
$ws.Range("C72").Value = 3.659
$ws.Range("D72").Value = 8.1155
$ws.Range("E72").Value = 27.1239
$ws.Range("B73").Value = "RHIM"
$ws.Range("C73").Value = 3.6544
$ws.Range("D73").Value = 3.2276
$ws.Range("E73").Value = 5.1826
$ws.Range("B74").Value = "CGPOWER"
$ws.Range("C74").Value = 3.6125
$ws.Range("D74").Value = 3.4192
$ws.Range("E74").Value = 1.0325
$ws.Range("B75").Value = "WELSPUNLIV"
$ws.Range("C75").Value = 3.6073
$ws.Range("D75").Value = 3.7285
$ws.Range("E75").Value = 15.9372
$ws.Range("B76").Value = "PFC"
$ws.Range("C76").Value = 3.5986
$ws.Range("D76").Value = 3.8618
$ws.Range("E76").Value = -0.3534

$ws = $wb.Worksheets.Item("Top Losers")
$ws.Range("B15").Value = "CCCL"
$ws.Range("C15").Value = -5.0146
$ws.Range("D15").Value = -4.576
$ws.Range("E15").Value = -12.4759
$ws.Range("B16").Value = "CREDITACC"
$ws.Range("C16").Value = -4.9692
$ws.Range("D16").Value = -1.3216
$ws.Range("E16").Value = 3.7319
$ws.Range("B17").Value = "KALAMANDIR"
$ws.Range("C17").Value = -4.8415
$ws.Range("D17").Value = 1.7451
$ws.Range("E17").Value = 25.9996
$ws.Range("B18").Value = "CRAMC"
$ws.Range("C18").Value = -4.7668
$ws.Range("D18").Value = 5.978
$ws.Range("E18").Value = "N/A"
$ws.Range("B19").Value = "SMLISUZU"
$ws.Range("C19").Value = -4.7654
$ws.Range("D19").Value = 4.993
$ws.Range("E19").Value = -2.8236
$ws.Range("B37").Value = "ANANDRATHI"
$ws.Range("C37").Value = -3.0775
$ws.Range("D37").Value = -0.8672
$ws.Range("E37").Value = 9.1835
$ws.Range("B38").Value = "NLCINDIA"
$ws.Range("C38").Value = -3.0757
$ws.Range("D38").Value = -4.5618
$ws.Range("E38").Value = -11.6431
$ws.Range("B39").Value = "YATRA"
$ws.Range("C39").Value = -3.0403
$ws.Range("D39").Value = -2.8455
$ws.Range("E39").Value = 7.3711
$ws.Range("B40").Value = "MPSLTD"
$ws.Range("C40").Value = -3.0335
$ws.Range("D40").Value = -4.3902
$ws.Range("E40").Value = 2.434
$ws.Range("B41").Value = "DRREDDY"
$ws.Range("C41").Value = -2.9859
$ws.Range("D41").Value = -2.5475
$ws.Range("E41").Value = 2.2228
$ws.Range("B42").Value = "ROSSTECH"
$ws.Range("C42").Value = -2.9778
$ws.Range("D42").Value = 1.9028
$ws.Range("E42").Value = -6.8057
$ws.Range("B43").Value = "OAL"
$ws.Range("C43").Value = -2.9496
$ws.Range("D43").Value = -1.278
$ws.Range("E43").Value = 8.7362
$ws.Range("B44").Value = "ENDURANCE"
$ws.Range("C44").Value = -2.939
$ws.Range("D44").Value = -2.2945
$ws.Range("E44").Value = 3.4531
$ws.Range("B45").Value = "POLICYBZR"
$ws.Range("C45").Value = -2.907
$ws.Range("D45").Value = 2.2365
$ws.Range("E45").Value = 1.2573
$ws.Range("B46").Value = "BOSCHLTD"
$ws.Range("C46").Value = -2.9061
$ws.Range("D46").Value = -3.0193
$ws.Range("E46").Value = -1.9006
$ws.Range("B47").Value = "DIGITIDE"
$ws.Range("C47").Value = -2.8795
$ws.Range("D47").Value = 3.2317
$ws.Range("E47").Value = 6.2968
$ws.Range("B48").Value = "RUBICON"
$ws.Range("C48").Value = -2.8712
$ws.Range("D48").Value = -2.9654
$ws.Range("E48").Value = "N/A"
$ws.Range("B49").Value = "STARHEALTH"
$ws.Range("C49").Value = -2.8707
$ws.Range("D49").Value = -1.5572
$ws.Range("E49").Value = 7.5434
$ws.Range("B50").Value = "KIRIINDUS"
$ws.Range("C50").Value = -2.8411
$ws.Range("D50").Value = -1.3849
$ws.Range("E50").Value = 1.4335
$ws.Range("B51").Value = "TTKPRESTIG"
$ws.Range("C51").Value = -2.7438
$ws.Range("D51").Value = 8.0012
$ws.Range("E51").Value = 9.6505
$ws.Range("B52").Value = "PFOCUS"
$ws.Range("C52").Value = -2.7039
$ws.Range("D52").Value = -2.6276
$ws.Range("E52").Value = -1.2163

$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Range("C12").Value = 45.3321
$ws.Range("C19").Value = 36.4813
$ws.Range("B71").Value = "THOMASCOTT"
$ws.Range("C71").Value = 19.1649
$ws.Range("B72").Value = "KARURVYSYA"
$ws.Range("C72").Value = 19.11
$ws.Range("B73").Value = "IIFL"
$ws.Range("C73").Value = 18.9853
$ws.Range("B74").Value = "LUMAXIND"
$ws.Range("C74").Value = 18.9608
